# S_Korea_scaling_mapping.xlsx update
# Remove shipping (Ship / value 18) and aviation (Domestic-aviation / value 28)
# duplicate entries from column C (and the stray B-column duplicates on the
# "2B_Chemical-industry"/"2D_Degreasing-Cleaning" header rows), leaving the
# single canonical mapping per ceds_sector row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("map")
$ws.Activate()

# Row 18 ("Aviation" / 1A3aii_Domestic-aviation row): drop the extra B18
# duplicate, keep C18.
$ws.Range("B18").ClearContents()

# Rows 20-26 (Road / 1A3b_Road block): drop the extra C column duplicates,
# keep column B.
$ws.Range("C20").ClearContents()
$ws.Range("C21").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Range("C23").ClearContents()
$ws.Range("C24").ClearContents()
$ws.Range("C25").ClearContents()
$ws.Range("C26").ClearContents()

# Row 28 ("Ship" / 2B_Chemical-industry row): drop the extra B28 duplicate,
# keep C28.
$ws.Range("B28").ClearContents()

# Row 32 (Construction facilities row): drop the stray C32 value that no
# longer matches B32.
$ws.Range("C32").ClearContents()

# Row 43 (Inorganic chemical manufacturing): drop the duplicate C43.
$ws.Range("C43").ClearContents()

# Row 51 (2D_Degreasing-Cleaning row): drop the duplicate C51.
$ws.Range("C51").ClearContents()

# Row 56 (Food and beverage processing row): drop the duplicate C56.
$ws.Range("C56").ClearContents()

# Update the view state: the frozen pane now shows row 41 at the top, and the
# last active cell in the lower pane is C73.
$ws.Range("C73").Select()
